$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autonomous_temporary")

# --- Fix typos in the "Choices" column for the Mission_Select signals ---
$ws.Range("K3").Value = "0=Manual, 1=Acceleration, 2=Skidpad, 3=Trackdrive, 4=Braketest, 5=Inspection, 6=Autocross"
$ws.Range("K7").Value = "0=Manual, 1=Acceleation, 2=Skidpadd, 3=Trackdive, 4=Braketest, 5=Inspection, 6=Autocross"

# --- Insert a new "Emergency" signal row under the ACU_IGN message (row 21), ---
# --- pushing everything from the old row 21 onward down by one row.          ---
$ws.Rows.Item(21).Insert()

# Copy the formatting of the row above (a data row) into the new row A21:K21
# so the new row matches the look of the other signal rows.
$ws.Range("A20:K20").Copy()
$ws.Range("A21:K21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A21").Value = "Emergency"
$ws.Range("B21").Value = 16
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = "Intel"
$ws.Range("E21").Value = $false
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0
$ws.Range("K21").Value = "1=Emergency_ON, 0=Emergency_OFF"
